$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.896.14"
$ws.Range("E2").Value = "  +0.01%  "

$ws.Range("D3").Value = "2.918.42"
$ws.Range("E3").Value = "  -0.11%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.66"
$ws.Range("E5").Value = "  +1.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.24"
$ws.Range("E6").Value = "  +1.03%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.507"
$ws.Range("E8").Value = "  +0.52%  "

$ws.Range("E9").Value = "  +1.93%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.144"
$ws.Range("E10").Value = "  -0.79%  "

$ws.Range("E11").Value = "  -1.77%  "

$ws.Range("E12").Value = "  -0.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.58"
$ws.Range("E13").Value = "  -0.03%  "

$ws.Range("E14").Value = "  +0.02%  "

$ws.Range("D15").Value = "3.399.22"
$ws.Range("E15").Value = "  -0.17%  "

$ws.Range("D16").Value = "60.811.87"
$ws.Range("E16").Value = "  -0.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.69"
$ws.Range("E17").Value = "  -1.04%  "

$ws.Range("D18").Value = "2.916.61"
$ws.Range("E18").Value = "  -0.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "430.97"
$ws.Range("E19").Value = "  -0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.36"
$ws.Range("E20").Value = "  -2.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.679"
$ws.Range("E21").Value = "  -0.80%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.05"
$ws.Range("E22").Value = "  -1.44%  "

$ws.Range("E23").Value = "  +1.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.04"
$ws.Range("E24").Value = "  +1.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.19"
$ws.Range("E25").Value = "  -0.89%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.82"
$ws.Range("E26").Value = "  -0.60%  "

$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.28"
$ws.Range("E28").Value = "  +5.95%  "

$ws.Range("E29").Value = "  -0.25%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.04"
$ws.Range("E30").Value = "  -2.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.60"
$ws.Range("E31").Value = "  -0.03%  "

$ws.Range("E32").Value = "  +1.40%  "

$ws.Range("E33").Value = "  +0.10%  "

$ws.Range("D34").Value = "0.0₃0866"
$ws.Range("E34").Value = "  -1.11%  "

$ws.Range("E35").Value = "  +0.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.63"
$ws.Range("E36").Value = "  -0.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.02"
$ws.Range("E37").Value = "  +0.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.98"
$ws.Range("E38").Value = "  -1.59%  "

$ws.Range("E39").Value = "  -4.42%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.56"
$ws.Range("E40").Value = "  -1.76%  "

$ws.Range("E41").Value = "  -3.41%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.40"
$ws.Range("E42").Value = "  -3.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "380.39"
$ws.Range("E43").Value = "  +0.32%  "

$ws.Range("D44").Value = "2.694.96"
$ws.Range("E44").Value = "  +0.59%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0343"
$ws.Range("E45").Value = "  -1.80%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.42"
$ws.Range("E46").Value = "  +0.62%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.86"
$ws.Range("E48").Value = "  -2.50%  "

$ws.Range("E49").Value = "  -0.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.01"
$ws.Range("E50").Value = "  -2.89%  "

$ws.Range("E51").Value = "  +0.18%  "
